# Apply MassWateR Results Template update:
# - Rename the "Grab" sample collection method ID to "Grab-MassWateR"
#   (and update the related instructions text on the Instructions sheet
#   accordingly) throughout the workbook.
# - Update the active-cell selections left on the Results and Instructions
#   sheets.

$wb = $excel.ActiveWorkbook
$wsResults = $wb.Worksheets.Item("Results")
$wsInstructions = $wb.Worksheets.Item("Instructions")

# --- Results sheet: example row used "Grab" as the Sample Collection Method ID
$wsResults.Range("O3").Value = "Grab-MassWateR"

# --- Instructions sheet: Sample Collection Method ID row (row 20)
$wsInstructions.Range("C20").Value = "Grab-MassWateR"
$wsInstructions.Range("B20").Value = "For WQX:  Enter the method ID used for this sample collection.  Not applicable for field measurement/observations.  Method IDs are defined in WQX by organization.  MassWateR will assign a default value of ""Grab-MassWateR"" if nothing is entered, but this requires a Method Context of ""MassWateR"" in the WQXMeta file.  Standard method IDs that can be used by any organization under the MassWateR context are ""Grab-MassWateR"", ""Pole-MassWateR"", and ""Basket-MassWateR""."

# --- Update leftover active-cell selections on each sheet
$wsInstructions.Activate()
$wsInstructions.Range("A6").Select() | Out-Null

$wsResults.Activate()
$wsResults.Range("A3").Select() | Out-Null
